# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview!E2:F3, and the "Status" column (C) on the
#    per-locale sheets zh-cn / de-de).
# 2. Narrow the "Status" column(s) (Overview!E:F and Status col on the
#    locale sheets) to their new, narrower width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Text replacement -----------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Column width changes --------------------------------------------
# Target stored (XML) column width is ~13.41 characters; this runtime
# quantizes ColumnWidth to a whole-pixel grid, so 12.5 (mid-bucket) is
# the value that lands on the nearest representable width.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
